$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Relabel the repeated column headers with numeric suffixes (1/2/3) for each
# of the three SATUAN/HARGABELI/HARGAJUAL/JMLBARANG/BARCODE blocks.
$ws.Range("G1").Value = "SATUAN1"
$ws.Range("H1").Value = "HARGABELI1"
$ws.Range("I1").Value = "HARGAJUAL1"
$ws.Range("J1").Value = "JMLBARANG1"
$ws.Range("K1").Value = "BARCODE1"

$ws.Range("M1").Value = "SATUAN2"
$ws.Range("N1").Value = "HARGABELI2"
$ws.Range("O1").Value = "HARGAJUAL2"
$ws.Range("P1").Value = "JMLBARANG2"
$ws.Range("Q1").Value = "BARCODE2"

$ws.Range("S1").Value = "SATUAN3"
$ws.Range("T1").Value = "HARGABELI3"
$ws.Range("U1").Value = "HARGAJUAL3"
$ws.Range("V1").Value = "JMLBARANG3"
$ws.Range("W1").Value = "BARCODE3"

# The three repeated blocks (SATUAN/HARGABELI/HARGAJUAL/JMLBARANG/BARCODE)
# now all share identical header text, so Excel's "best fit" column widths
# recompute to the same values across each of the 3 blocks.
$ws.Range("G1").ColumnWidth = 8.5
$ws.Range("H1").ColumnWidth = 11
$ws.Range("I1").ColumnWidth = 11.666666666666666
$ws.Range("J1").ColumnWidth = 12
$ws.Range("K1").ColumnWidth = 9.666666666666666

$ws.Range("M1").ColumnWidth = 8.5
$ws.Range("N1").ColumnWidth = 11
$ws.Range("O1").ColumnWidth = 11.666666666666666
$ws.Range("P1").ColumnWidth = 12

$ws.Range("S1").ColumnWidth = 8.5
$ws.Range("T1").ColumnWidth = 11
$ws.Range("U1").ColumnWidth = 11.666666666666666
$ws.Range("V1").ColumnWidth = 12
$ws.Range("W1").ColumnWidth = 9.666666666666666

# Update the view: scroll so column E is the left-most visible column and
# move the active selection to T10.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("T10").Select()
